$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - DataHub & MUMC+ (years 2019-2023 -> columns G,H,I,J,K)
$ws.Range("G2").Value = 10
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 20
$ws.Range("J2").Value = 20
$ws.Range("K2").Value = 40

# Row 3 - FASoS (years 2020-2023 -> columns H,I,J,K)
$ws.Range("H3").Value = 16.67
$ws.Range("I3").Value = 33.33
$ws.Range("J3").Value = 16.67
$ws.Range("K3").Value = 33.33

# Row 4 - FHML (years 2015-2018, 2020-2023 -> columns C,D,E,F,H,I,J,K)
$ws.Range("C4").Value = 3.17
$ws.Range("D4").Value = 6.35
$ws.Range("E4").Value = 9.52
$ws.Range("F4").Value = 11.11
$ws.Range("H4").Value = 14.29
$ws.Range("I4").Value = 17.46
$ws.Range("J4").Value = 26.98
$ws.Range("K4").Value = 11.11

# Row 5 - FPN (years 2014-2023 -> columns B,C,D,E,F,G,H,I,J,K)
$ws.Range("B5").Value = 4.57
$ws.Range("C5").Value = 1.02
$ws.Range("D5").Value = 1.52
$ws.Range("E5").Value = 5.08
$ws.Range("F5").Value = 6.09
$ws.Range("G5").Value = 13.2
$ws.Range("H5").Value = 18.78
$ws.Range("I5").Value = 22.34
$ws.Range("J5").Value = 18.27
$ws.Range("K5").Value = 9.140000000000001
